$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.3033442088659019
$ws.Range("C2").Value = -0.285613861520414
$ws.Range("B3").Value = 0.2207520964114534
$ws.Range("C3").Value = 0.2309543543290238
$ws.Range("B4").Value = 0.1631585845063193
$ws.Range("C4").Value = 0.1757923210546981
$ws.Range("B5").Value = -0.05055669449508423
$ws.Range("C5").Value = -0.04104158479270387
$ws.Range("B6").Value = 0.1569216718066155
$ws.Range("C6").Value = 0.1728512757315558
$ws.Range("B7").Value = -0.4385607354708777
$ws.Range("C7").Value = -0.4281640872318672
$ws.Range("B8").Value = -0.2815207651310452
$ws.Range("C8").Value = -0.2597637584727662
$ws.Range("B9").Value = -0.4053244950965031
$ws.Range("C9").Value = -0.3948693068931027
$ws.Range("B10").Value = 0.3291626202766235
$ws.Range("C10").Value = 0.3465381306558325
$ws.Range("B11").Value = -0.2026559651893725
$ws.Range("C11").Value = -0.1947378753737388
$ws.Range("B12").Value = -0.08357465415811693
$ws.Range("C12").Value = -0.05811555757293763
$ws.Range("B13").Value = 0.003268724058882566
$ws.Range("C13").Value = 0.003907374428422664
$ws.Range("B14").Value = 0.03850792956741418
$ws.Range("C14").Value = 0.07863851839863101
$ws.Range("B15").Value = -0.07913911116193134
$ws.Range("C15").Value = -0.02988053198258833
$ws.Range("B16").Value = 0.2662829113388356
$ws.Range("C16").Value = 0.3323814523503869
$ws.Range("B17").Value = 0.5778272802240982
$ws.Range("C17").Value = 0.6103798634913413
$ws.Range("B18").Value = 0.05827459139698862
$ws.Range("C18").Value = 0.04020067918106943
$ws.Range("B19").Value = 0.3999590293933805
$ws.Range("C19").Value = 0.4091454085077181
$ws.Range("B20").Value = 0.2406752051786359
$ws.Range("C20").Value = 0.2999452911663952
$ws.Range("B21").Value = 0.4417628114591178
$ws.Range("C21").Value = 0.5085908455124964
$ws.Range("B22").Value = 0.3494371019363526
$ws.Range("C22").Value = 0.3861123216915772
$ws.Range("B23").Value = -0.06659496235699991
$ws.Range("C23").Value = -0.03456165225398219
$ws.Range("B24").Value = 4.541078537018452
$ws.Range("C24").Value = 4.566046494845523
$ws.Range("B25").Value = 0.5449400901430083
$ws.Range("C25").Value = 0.5104499532924316
$ws.Range("B26").Value = 0.4002625929779248
$ws.Range("C26").Value = 0.3821785279846542
$ws.Range("B27").Value = 0.3312475927322071
$ws.Range("C27").Value = 0.3022487667381129
$ws.Range("B28").Value = 1.103864868107481
$ws.Range("C28").Value = 1.073524239157886
$ws.Range("B29").Value = 5.837247815210643
$ws.Range("C29").Value = 5.385828885282493
$ws.Range("B30").Value = 1.015620096334174
$ws.Range("C30").Value = 0.9650930574251766
$ws.Range("B31").Value = -0.128165917865441
$ws.Range("C31").Value = -0.1933917469891227
$ws.Range("B32").Value = 0.8107004350220448
$ws.Range("C32").Value = 0.7727241770431046
$ws.Range("B33").Value = 0.9162711559932648
$ws.Range("C33").Value = 0.8903333690452184
$ws.Range("B34").Value = -0.5988286035607752
$ws.Range("C34").Value = -0.6291079321617313
$ws.Range("B35").Value = 0.8319861661415534
$ws.Range("C35").Value = 0.8220208498997641
$ws.Range("B36").Value = 0.7854255156702361
$ws.Range("C36").Value = 0.7692466651518581
$ws.Range("B37").Value = 0.7685533230710442
$ws.Range("C37").Value = 0.7478811412401991
$ws.Range("B38").Value = 0.7594994496718845
$ws.Range("C38").Value = 0.7366819587082833
$ws.Range("B39").Value = 0.5831427592561873
$ws.Range("C39").Value = 0.5807328869148602
$ws.Range("B40").Value = 0.7542173216697582
$ws.Range("C40").Value = 0.7527880614441489
$ws.Range("B41").Value = 0.5740814366689742
$ws.Range("C41").Value = 0.5663686018189583
$ws.Range("B42").Value = 0.7219259599136361
$ws.Range("C42").Value = 0.6924706582687463
$ws.Range("B43").Value = 0.733890551282992
$ws.Range("C43").Value = 0.7182193194617948
$ws.Range("B44").Value = 0.6745894397723444
$ws.Range("C44").Value = 0.6674199276345331
$ws.Range("B45").Value = 0.6604177017042322
$ws.Range("C45").Value = 0.643203384512019
$ws.Range("B46").Value = -1.251814518125235
$ws.Range("C46").Value = -1.256094731000866
$ws.Range("B47").Value = -0.9693432235046766
$ws.Range("C47").Value = -0.9744620795456761
$ws.Range("B48").Value = -0.862326360188199
$ws.Range("C48").Value = -0.8685574750665104
$ws.Range("B49").Value = -0.6314247653529417
$ws.Range("C49").Value = -0.6349189543844442
$ws.Range("B50").Value = -0.04861888039615835
$ws.Range("C50").Value = -0.0490296225981823
$ws.Range("B51").Value = -0.8510358802390722
$ws.Range("C51").Value = -0.8550403088152924
$ws.Range("B52").Value = -0.8510358802390722
$ws.Range("C52").Value = -0.8550403088152924
$ws.Range("B53").Value = -1.078588183733273
$ws.Range("C53").Value = -1.093161550357953
$ws.Range("B54").Value = -0.1865304266581551
$ws.Range("C54").Value = -0.1850210583711067
$ws.Range("B55").Value = -0.9896890869033197
$ws.Range("C55").Value = -0.9934037464819662
$ws.Range("B56").Value = -0.8915757557210402
$ws.Range("C56").Value = -0.8858407144749817
$ws.Range("B57").Value = -0.9654658271717568
$ws.Range("C57").Value = -0.9484025178377414
$ws.Range("B58").Value = -1.166399175072524
$ws.Range("C58").Value = -1.135877974028292
$ws.Range("B59").Value = -0.8697146221930475
$ws.Range("C59").Value = -0.8520611867384603
$ws.Range("B60").Value = -0.5226236749088476
$ws.Range("C60").Value = -0.4967963367636913
$ws.Range("B61").Value = 0.3667234597952868
$ws.Range("C61").Value = 0.3698122843417765
$ws.Range("B62").Value = -1.256585013455179
$ws.Range("C62").Value = -1.233472899209369
$ws.Range("B63").Value = -0.76328821103421
$ws.Range("C63").Value = -0.7279101236544864
$ws.Range("B64").Value = -0.9019819226868909
$ws.Range("C64").Value = -0.8933700897086743
$ws.Range("B65").Value = -0.1414601443187503
$ws.Range("C65").Value = -0.1159650726353796
$ws.Range("B66").Value = -0.8343463378780934
$ws.Range("C66").Value = -0.8036033375799481
$ws.Range("B67").Value = -0.8369362897989698
$ws.Range("C67").Value = -0.7947377314035913

Write-Host "applied updates"
